# Update gh-pages output (regenerate "想去人数" / want-to-go counts) to match
# the run generated at 456a3b4, and append the newly-scraped local-life event
# ("本地生活" sheet) that the crawl picked up since the previous run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions) — column F ("想去人数") bumps
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 6456
$ws1.Range("F3").Value  = 2608
$ws1.Range("F5").Value  = 1312
$ws1.Range("F7").Value  = 3190
$ws1.Range("F8").Value  = 380
$ws1.Range("F11").Value = 8098
$ws1.Range("F12").Value = 420
$ws1.Range("F13").Value = 75
$ws1.Range("F16").Value = 289
$ws1.Range("F18").Value = 56
$ws1.Range("F19").Value = 497
$ws1.Range("F20").Value = 310
$ws1.Range("F21").Value = 10081
$ws1.Range("F26").Value = 378
$ws1.Range("F28").Value = 24
$ws1.Range("F29").Value = 139
$ws1.Range("F30").Value = 81
$ws1.Range("F33").Value = 2059
$ws1.Range("F34").Value = 28
$ws1.Range("F35").Value = 29
$ws1.Range("F37").Value = 4026
$ws1.Range("F38").Value = 252
$ws1.Range("F40").Value = 2092
$ws1.Range("F41").Value = 1214
$ws1.Range("F42").Value = 139
$ws1.Range("F43").Value = 296
$ws1.Range("F44").Value = 209
$ws1.Range("F46").Value = 87
$ws1.Range("F47").Value = 76
$ws1.Range("F49").Value = 51

# ---------------------------------------------------------------------
# Sheet "演出" (Performances) — column F ("想去人数") bumps
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value  = 8
$ws2.Range("F6").Value  = 36
$ws2.Range("F9").Value  = 26
$ws2.Range("F16").Value = 171
$ws2.Range("F18").Value = 23
$ws2.Range("F20").Value = 20

# ---------------------------------------------------------------------
# Sheet "本地生活" (Local life) — newly scraped event appended as row 2
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

# Column A carries the same bold/centered/bordered style as the header's
# row-index column; copy it from A1 so A2 matches the rest of the sheets'
# numbered rows (A1.Style is a COM object — value-copy, not formula, so a
# copy/paste-formats round trip is the reliable way to clone it).
$ws3.Range("A2").Value = 1
$ws3.Range("A1").Copy()
$ws3.Range("A2").PasteSpecial(-4122)

# "2024-05-03" looks like a date to Excel's literal-value parser, so force
# text formatting before assigning it (matching every other date-like "B"
# column cell in this workbook, which are all stored as plain text), then
# drop back to the default style so no stray number format lingers.
$ws3.Range("B2").NumberFormat = "@"
$ws3.Range("B2").Value = "2024-05-03"
$ws3.Range("B2").Style = "Normal"

$ws3.Range("C2").Value = "北京·塔罗集市"
$ws3.Range("D2").Value = "北京朝阳区广渠路南侧汇泰大厦1层展厅 汇泰大厦"
$ws3.Range("E2").Value = "2024.05.03 09:30-05.03 16:30"
$ws3.Range("F2").Value = 1
$ws3.Range("G2").Value = 56
$ws3.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=84239"
$ws3.Range("I2").Value = "//i2.hdslb.com/bfs/openplatform/202404/z8qnfmoq1712735872200.jpeg"

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types) — column F ("想去人数") bumps
# (row numbers differ slightly from "展览" because this combined sheet
# interleaves exhibitions + performances)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6456
$ws4.Range("F3").Value  = 2608
$ws4.Range("F6").Value  = 1312
$ws4.Range("F8").Value  = 3190
$ws4.Range("F9").Value  = 380
$ws4.Range("F13").Value = 8098
$ws4.Range("F14").Value = 420
$ws4.Range("F15").Value = 75
$ws4.Range("F18").Value = 289
$ws4.Range("F19").Value = 56
$ws4.Range("F20").Value = 497
$ws4.Range("F21").Value = 310
$ws4.Range("F22").Value = 10081
$ws4.Range("F26").Value = 378
$ws4.Range("F29").Value = 139
$ws4.Range("F30").Value = 81
$ws4.Range("F33").Value = 2059
$ws4.Range("F34").Value = 28
$ws4.Range("F36").Value = 4026
$ws4.Range("F37").Value = 252
$ws4.Range("F39").Value = 2093
$ws4.Range("F40").Value = 23
$ws4.Range("F41").Value = 1214
$ws4.Range("F42").Value = 139
$ws4.Range("F43").Value = 296
$ws4.Range("F44").Value = 209
$ws4.Range("F46").Value = 87
$ws4.Range("F47").Value = 76
$ws4.Range("F49").Value = 51
